$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 1.187716
$ws.Range("H2").Value = 3.563148
$ws.Range("M2").Value = 1.837384
$ws.Range("N2").Value = 5.512152
$ws.Range("O2").Value = 0.0635335947613339
$ws.Range("P2").Value = 0.0635335947613339
$ws.Range("Q2").Value = 2.182290374944
$ws.Range("R2").Value = 19.640613374496
$ws.Range("S2").Value = 0.0635335947613339
$ws.Range("T2").Value = 0.0635335947613339

# Row 3
$ws.Range("G3").Value = 1.187716
$ws.Range("H3").Value = 3.563148
$ws.Range("O3").Value = 0.02082867030699976
$ws.Range("P3").Value = 0.02082867030699976
$ws.Range("Q3").Value = 0.7154357770026667
$ws.Range("R3").Value = 6.438921993024
$ws.Range("S3").Value = 0.02082867030699976
$ws.Range("T3").Value = 0.02082867030699976

# Row 4
$ws.Range("G4").Value = 1.187716
$ws.Range("H4").Value = 3.563148
$ws.Range("M4").Value = 5.528959666666668
$ws.Range("N4").Value = 16.586879
$ws.Range("O4").Value = 0.1911819646376369
$ws.Range("P4").Value = 0.1911819646376369
$ws.Range("Q4").Value = 6.566833859454669
$ws.Range("R4").Value = 59.10150473509202
$ws.Range("S4").Value = 0.1911819646376369
$ws.Range("T4").Value = 0.1911819646376369

# Row 5
$ws.Range("G5").Value = 1.187716
$ws.Range("H5").Value = 3.563148
$ws.Range("M5").Value = 3.046454666666667
$ws.Range("N5").Value = 9.139364
$ws.Range("O5").Value = 0.1053411895666744
$ws.Range("P5").Value = 0.1053411895666744
$ws.Range("Q5").Value = 3.618322950874667
$ws.Range("R5").Value = 32.56490655787201
$ws.Range("S5").Value = 0.1053411895666744
$ws.Range("T5").Value = 0.1053411895666744

# Row 6
$ws.Range("G6").Value = 1.187716
$ws.Range("H6").Value = 3.563148
$ws.Range("M6").Value = 9.021246333333332
$ws.Range("N6").Value = 27.063739
$ws.Range("O6").Value = 0.3119392618985303
$ws.Range("P6").Value = 0.3119392618985303
$ws.Range("Q6").Value = 10.71467861004133
$ws.Range("R6").Value = 96.432107490372
$ws.Range("S6").Value = 0.3119392618985303
$ws.Range("T6").Value = 0.3119392618985303

# Row 7
$ws.Range("G7").Value = 1.187716
$ws.Range("H7").Value = 3.563148
$ws.Range("M7").Value = 8.883473666666667
$ws.Range("N7").Value = 26.650421
$ws.Range("O7").Value = 0.3071753188288246
$ws.Range("P7").Value = 0.3071753188288246
$ws.Range("Q7").Value = 10.55104380947867
$ws.Range("R7").Value = 94.95939428530802
$ws.Range("S7").Value = 0.3071753188288246
$ws.Range("T7").Value = 0.3071753188288246
